# Auto-generated edit script: updates Kraken_Profits price/profit figures
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(4, 8).Value = 2086.0908  # H4
$ws.Cells.Item(4, 9).Value = 2086.0908  # I4
$ws.Cells.Item(4, 11).Value = 2086.0908  # K4
$ws.Cells.Item(4, 13).Value = -1972.0908  # M4

$ws.Cells.Item(8, 8).Value = 9  # H8
$ws.Cells.Item(8, 9).Value = 9  # I8
$ws.Cells.Item(8, 11).Value = 27  # K8
$ws.Cells.Item(8, 13).Value = 112  # M8

$ws.Cells.Item(32, 8).Value = 8428.571  # H32
$ws.Cells.Item(32, 9).Value = 0  # I32
$ws.Cells.Item(32, 10).Value = 8428.571  # J32
$ws.Cells.Item(32, 11).Value = 0  # K32
$ws.Cells.Item(32, 12).Value = 8428.571  # L32
$ws.Cells.Item(32, 13).ClearContents()  # M32
$ws.Cells.Item(32, 14).Value = -9080.571  # N32

$ws.Cells.Item(64, 8).Value = 3950  # H64
$ws.Cells.Item(64, 9).Value = 4000  # I64
$ws.Cells.Item(64, 10).Value = 3900  # J64
$ws.Cells.Item(64, 11).Value = 4000  # K64
$ws.Cells.Item(64, 12).Value = 3900  # L64
$ws.Cells.Item(64, 13).Value = -3752  # M64
$ws.Cells.Item(64, 14).Value = -4396  # N64

$ws.Cells.Item(67, 8).Value = 3950  # H67
$ws.Cells.Item(67, 9).Value = 4000  # I67
$ws.Cells.Item(67, 10).Value = 3900  # J67
$ws.Cells.Item(67, 11).Value = 4000  # K67
$ws.Cells.Item(67, 12).Value = 3900  # L67
$ws.Cells.Item(67, 13).Value = -3142  # M67
$ws.Cells.Item(67, 14).Value = -5616  # N67

$ws.Cells.Item(70, 8).Value = 1575  # H70
$ws.Cells.Item(70, 10).Value = 2150  # J70
$ws.Cells.Item(70, 12).Value = 6450  # L70
$ws.Cells.Item(70, 14).Value = -6990  # N70

$ws.Cells.Item(73, 8).Value = 1575  # H73
$ws.Cells.Item(73, 10).Value = 2150  # J73
$ws.Cells.Item(73, 12).Value = 6450  # L73
$ws.Cells.Item(73, 14).Value = -8322  # N73

$ws.Cells.Item(86, 8).Value = 2000  # H86
$ws.Cells.Item(86, 10).Value = 0  # J86
$ws.Cells.Item(86, 12).Value = 0  # L86
$ws.Cells.Item(86, 14).ClearContents()  # N86

$ws.Cells.Item(89, 8).Value = 2000  # H89
$ws.Cells.Item(89, 10).Value = 0  # J89
$ws.Cells.Item(89, 12).Value = 0  # L89
$ws.Cells.Item(89, 14).ClearContents()  # N89

$ws.Cells.Item(107, 8).Value = 3099.8  # H107
$ws.Cells.Item(107, 9).Value = 2333.3333  # I107
$ws.Cells.Item(107, 11).Value = 2333.3333  # K107
$ws.Cells.Item(107, 13).Value = -413.3332999999998  # M107

$ws.Cells.Item(116, 8).Value = 4098.3335  # H116
$ws.Cells.Item(116, 9).Value = 3773.75  # I116
$ws.Cells.Item(116, 11).Value = 3773.75  # K116
$ws.Cells.Item(116, 13).Value = -331.75  # M116

$ws.Cells.Item(132, 8).Value = 4326.4287  # H132
$ws.Cells.Item(132, 9).Value = 2955.9  # I132
$ws.Cells.Item(132, 10).Value = 7752.75  # J132
$ws.Cells.Item(132, 11).Value = 8867.700000000001  # K132
$ws.Cells.Item(132, 12).Value = 23258.25  # L132
$ws.Cells.Item(132, 13).Value = -6337.700000000001  # M132
$ws.Cells.Item(132, 14).Value = -28318.25  # N132

$ws.Cells.Item(138, 8).Value = 3075.16  # H138
$ws.Cells.Item(138, 9).Value = 1112  # I138
$ws.Cells.Item(138, 10).Value = 3999  # J138
$ws.Cells.Item(138, 11).Value = 3336  # K138
$ws.Cells.Item(138, 12).Value = 11997  # L138
$ws.Cells.Item(138, 13).Value = 1804  # M138
$ws.Cells.Item(138, 14).Value = -22277  # N138

$ws.Cells.Item(141, 8).Value = 5895.5  # H141
$ws.Cells.Item(141, 9).Value = 5895.5  # I141
$ws.Cells.Item(141, 11).Value = 17686.5  # K141
$ws.Cells.Item(141, 13).Value = -12506.5  # M141

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(43, 8).Value = 44447.332  # H43
$ws.Cells.Item(43, 10).Value = 43500  # J43
$ws.Cells.Item(43, 12).Value = 43500  # L43
$ws.Cells.Item(43, 14).Value = -44126  # N43

$ws.Cells.Item(45, 8).Value = 2773.7144  # H45
$ws.Cells.Item(45, 9).Value = 2280.4  # I45
$ws.Cells.Item(45, 11).Value = 2280.4  # K45
$ws.Cells.Item(45, 13).Value = -1903.4  # M45

$ws.Cells.Item(53, 8).Value = 0  # H53
$ws.Cells.Item(53, 9).Value = 0  # I53
$ws.Cells.Item(53, 11).Value = 0  # K53
$ws.Cells.Item(53, 13).ClearContents()  # M53

$ws.Cells.Item(82, 8).Value = 0  # H82
$ws.Cells.Item(82, 10).Value = 0  # J82
$ws.Cells.Item(82, 12).Value = 0  # L82
$ws.Cells.Item(82, 14).ClearContents()  # N82

$ws.Cells.Item(85, 8).Value = 0  # H85
$ws.Cells.Item(85, 10).Value = 0  # J85
$ws.Cells.Item(85, 12).Value = 0  # L85
$ws.Cells.Item(85, 14).ClearContents()  # N85

$ws.Cells.Item(97, 8).Value = 3076.75  # H97
$ws.Cells.Item(97, 9).Value = 2060.8572  # I97
$ws.Cells.Item(97, 10).Value = 4499  # J97
$ws.Cells.Item(97, 11).Value = 2060.8572  # K97
$ws.Cells.Item(97, 12).Value = 4499  # L97
$ws.Cells.Item(97, 13).Value = -1564.8572  # M97
$ws.Cells.Item(97, 14).Value = -5491  # N97

$ws.Cells.Item(132, 8).Value = 1531.7142  # H132
$ws.Cells.Item(132, 9).Value = 1341.9231  # I132
$ws.Cells.Item(132, 11).Value = 4025.7693  # K132
$ws.Cells.Item(132, 13).Value = -1495.7693  # M132

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 5000  # H20
$ws.Cells.Item(20, 9).Value = 5000  # I20
$ws.Cells.Item(20, 11).Value = 5000  # K20
$ws.Cells.Item(20, 13).Value = -4753  # M20

$ws.Cells.Item(76, 8).Value = 21749.5  # H76
$ws.Cells.Item(76, 10).Value = 21749.5  # J76
$ws.Cells.Item(76, 12).Value = 21749.5  # L76
$ws.Cells.Item(76, 14).Value = -22379.5  # N76

$ws.Cells.Item(79, 8).Value = 21749.5  # H79
$ws.Cells.Item(79, 10).Value = 21749.5  # J79
$ws.Cells.Item(79, 12).Value = 21749.5  # L79
$ws.Cells.Item(79, 14).Value = -23933.5  # N79

$ws.Cells.Item(99, 8).Value = 3454.4546  # H99
$ws.Cells.Item(99, 9).Value = 3599.9  # I99
$ws.Cells.Item(99, 11).Value = 3599.9  # K99
$ws.Cells.Item(99, 13).Value = -2101.9  # M99

$ws.Cells.Item(134, 8).Value = 4310.9473  # H134
$ws.Cells.Item(134, 9).Value = 3484.6155  # I134
$ws.Cells.Item(134, 10).Value = 6101.3335  # J134
$ws.Cells.Item(134, 11).Value = 10453.8465  # K134
$ws.Cells.Item(134, 12).Value = 18304.0005  # L134
$ws.Cells.Item(134, 13).Value = -7918.8465  # M134
$ws.Cells.Item(134, 14).Value = -23374.0005  # N134

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 10889.223  # H31
$ws.Cells.Item(31, 9).Value = 12802.4  # I31
$ws.Cells.Item(31, 10).Value = 8497.75  # J31
$ws.Cells.Item(31, 11).Value = 12802.4  # K31
$ws.Cells.Item(31, 12).Value = 8497.75  # L31
$ws.Cells.Item(31, 13).Value = -12507.4  # M31
$ws.Cells.Item(31, 14).Value = -9087.75  # N31

$ws.Cells.Item(34, 8).Value = 10889.223  # H34
$ws.Cells.Item(34, 9).Value = 12802.4  # I34
$ws.Cells.Item(34, 10).Value = 8497.75  # J34
$ws.Cells.Item(34, 11).Value = 12802.4  # K34
$ws.Cells.Item(34, 12).Value = 8497.75  # L34
$ws.Cells.Item(34, 13).Value = -12600.4  # M34
$ws.Cells.Item(34, 14).Value = -8901.75  # N34

$ws.Cells.Item(132, 8).Value = 3237.3333  # H132
$ws.Cells.Item(132, 9).Value = 3237.3333  # I132
$ws.Cells.Item(132, 11).Value = 9711.999899999999  # K132
$ws.Cells.Item(132, 13).Value = -7181.999899999999  # M132

$ws.Cells.Item(134, 8).Value = 7000  # H134
$ws.Cells.Item(134, 9).Value = 0  # I134
$ws.Cells.Item(134, 10).Value = 7000  # J134
$ws.Cells.Item(134, 11).Value = 0  # K134
$ws.Cells.Item(134, 12).Value = 21000  # L134
$ws.Cells.Item(134, 13).ClearContents()  # M134
$ws.Cells.Item(134, 14).Value = -26070  # N134

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(23, 8).Value = 881.2  # H23
$ws.Cells.Item(23, 9).Value = 831.7143  # I23
$ws.Cells.Item(23, 10).Value = 924.5  # J23
$ws.Cells.Item(23, 11).Value = 2495.1429  # K23
$ws.Cells.Item(23, 12).Value = 2773.5  # L23
$ws.Cells.Item(23, 13).Value = -2260.1429  # M23
$ws.Cells.Item(23, 14).Value = -3243.5  # N23

$ws.Cells.Item(93, 8).Value = 9999  # H93
$ws.Cells.Item(93, 9).Value = 9999  # I93
$ws.Cells.Item(93, 11).Value = 29997  # K93
$ws.Cells.Item(93, 13).Value = -28125  # M93

$ws.Cells.Item(107, 8).Value = 249.42857  # H107
$ws.Cells.Item(107, 10).Value = 324.75  # J107
$ws.Cells.Item(107, 12).Value = 974.25  # L107
$ws.Cells.Item(107, 14).Value = -4814.25  # N107

$ws.Cells.Item(113, 10).Value = 888  # J113
$ws.Cells.Item(113, 12).Value = 2664  # L113
$ws.Cells.Item(113, 14).Value = -7004  # N113

$ws.Cells.Item(129, 8).Value = 4345  # H129
$ws.Cells.Item(129, 9).Value = 2460  # I129
$ws.Cells.Item(129, 10).Value = 10000  # J129
$ws.Cells.Item(129, 11).Value = 7380  # K129
$ws.Cells.Item(129, 12).Value = 30000  # L129
$ws.Cells.Item(129, 13).Value = -2380  # M129
$ws.Cells.Item(129, 14).Value = -40000  # N129

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(5, 8).Value = 2202.4  # H5

$ws.Cells.Item(48, 8).Value = 0  # H48
$ws.Cells.Item(48, 10).Value = 0  # J48
$ws.Cells.Item(48, 12).Value = 0  # L48
$ws.Cells.Item(48, 14).ClearContents()  # N48

$ws.Cells.Item(93, 8).Value = 90000  # H93
$ws.Cells.Item(93, 9).Value = 90000  # I93
$ws.Cells.Item(93, 11).Value = 90000  # K93
$ws.Cells.Item(93, 13).Value = -88128  # M93

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(45, 8).Value = 37000  # H45
$ws.Cells.Item(45, 9).Value = 37000  # I45
$ws.Cells.Item(45, 11).Value = 37000  # K45
$ws.Cells.Item(45, 13).Value = -36593  # M45

$ws.Cells.Item(46, 8).Value = 848.8  # H46
$ws.Cells.Item(46, 9).Value = 811  # I46
$ws.Cells.Item(46, 11).Value = 811  # K46
$ws.Cells.Item(46, 13).Value = -623  # M46

$ws.Cells.Item(55, 8).Value = 1121.5714  # H55
$ws.Cells.Item(55, 9).Value = 270.2  # I55
$ws.Cells.Item(55, 10).Value = 3250  # J55
$ws.Cells.Item(55, 11).Value = 270.2  # K55
$ws.Cells.Item(55, 12).Value = 3250  # L55
$ws.Cells.Item(55, 13).Value = -97.19999999999999  # M55
$ws.Cells.Item(55, 14).Value = -3596  # N55

$ws.Cells.Item(100, 8).Value = 6277.4287  # H100
$ws.Cells.Item(100, 10).Value = 12060.2  # J100
$ws.Cells.Item(100, 12).Value = 12060.2  # L100
$ws.Cells.Item(100, 14).Value = -13142.2  # N100

$ws.Cells.Item(136, 8).Value = 1237.25  # H136
$ws.Cells.Item(136, 9).Value = 1149.6666  # I136
$ws.Cells.Item(136, 11).Value = 3448.9998  # K136
$ws.Cells.Item(136, 13).Value = -898.9998000000001  # M136

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(22, 8).Value = 0  # H22
$ws.Cells.Item(22, 9).Value = 0  # I22
$ws.Cells.Item(22, 11).Value = 0  # K22
$ws.Cells.Item(22, 13).ClearContents()  # M22

$ws.Cells.Item(81, 8).Value = 1000  # H81
$ws.Cells.Item(81, 9).Value = 0  # I81
$ws.Cells.Item(81, 10).Value = 1000  # J81
$ws.Cells.Item(81, 11).Value = 0  # K81
$ws.Cells.Item(81, 12).Value = 2000  # L81
$ws.Cells.Item(81, 13).ClearContents()  # M81
$ws.Cells.Item(81, 14).Value = -4122  # N81

$ws.Cells.Item(84, 8).Value = 1000  # H84
$ws.Cells.Item(84, 9).Value = 0  # I84
$ws.Cells.Item(84, 10).Value = 1000  # J84
$ws.Cells.Item(84, 11).Value = 0  # K84
$ws.Cells.Item(84, 12).Value = 10000  # L84
$ws.Cells.Item(84, 13).ClearContents()  # M84
$ws.Cells.Item(84, 14).Value = -20608  # N84

$ws.Cells.Item(107, 8).Value = 350  # H107
$ws.Cells.Item(107, 9).Value = 350  # I107
$ws.Cells.Item(107, 11).Value = 1050  # K107
$ws.Cells.Item(107, 13).Value = 870  # M107

$ws.Cells.Item(122, 8).Value = 2625  # H122
$ws.Cells.Item(122, 9).Value = 2166.6667  # I122
$ws.Cells.Item(122, 10).Value = 4000  # J122
$ws.Cells.Item(122, 11).Value = 6500.000100000001  # K122
$ws.Cells.Item(122, 12).Value = 12000  # L122
$ws.Cells.Item(122, 13).Value = -4050.000100000001  # M122
$ws.Cells.Item(122, 14).Value = -16900  # N122

$ws.Cells.Item(132, 8).Value = 3266.2964  # H132
$ws.Cells.Item(132, 9).Value = 2202.8667  # I132
$ws.Cells.Item(132, 10).Value = 4595.5835  # J132
$ws.Cells.Item(132, 11).Value = 6608.6001  # K132
$ws.Cells.Item(132, 12).Value = 13786.7505  # L132
$ws.Cells.Item(132, 13).Value = -4078.6001  # M132
$ws.Cells.Item(132, 14).Value = -18846.7505  # N132
